$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active sheet (mirrors activeTab moving from "Edit Repayment
# Schedule" to "Repayment schedule" in the saved workbook view).
$ws.Activate()

# Insert a new (blank) column before column N, shifting the existing
# N/O/P ("Late", "heading"/Over Due, "Outstanding") columns one place right.
$ws.Columns("N").Insert()

# The newly inserted column picks up the width of its left neighbour (M).
$ws.Columns("N").ColumnWidth = 10.17

# Leave the cursor parked on S9, matching the saved selection.
$ws.Range("S9").Select()
